$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2392
$ws.Range("L3").Value = 2413
$ws.Range("L4").Value = 658
$ws.Range("H5").Value = 806
$ws.Range("L5").Value = 145
$ws.Range("K6").Value = 9122
$ws.Range("L6").Value = 2194
$ws.Range("H7").Value = 26069
$ws.Range("K7").Value = 27555
$ws.Range("L7").Value = 7802

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L4").Value = 39
$ws.Range("L7").Value = 496

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 96
$ws.Range("L6").Value = 122
$ws.Range("L7").Value = 355

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L6").Value = 88
$ws.Range("L7").Value = 286

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 29
$ws.Range("L6").Value = 60
$ws.Range("L7").Value = 252
$ws.Range("L8").Value = 496
$ws.Range("L11").Value = 137
$ws.Range("L18").Value = 57
$ws.Range("L19").Value = 222
$ws.Range("L20").Value = 197
$ws.Range("L23").Value = 81
$ws.Range("L29").Value = 398
$ws.Range("L33").Value = 355
$ws.Range("L37").Value = 286
$ws.Range("L40").Value = 20
$ws.Range("L42").Value = 247
$ws.Range("L44").Value = 58
$ws.Range("K48").Value = 348
$ws.Range("L51").Value = 88
$ws.Range("L52").Value = 154
$ws.Range("L57").Value = 34
$ws.Range("L61").Value = 10
$ws.Range("H63").Value = 305
$ws.Range("L66").Value = 19
$ws.Range("L67").Value = 287
$ws.Range("L76").Value = 90
$ws.Range("L79").Value = 216
$ws.Range("L84").Value = 77
$ws.Range("L85").Value = 411
$ws.Range("L91").Value = 112
$ws.Range("L93").Value = 41
$ws.Range("L94").Value = 92
$ws.Range("L99").Value = 123
$ws.Range("H101").Value = 26069
$ws.Range("K101").Value = 27555
$ws.Range("L101").Value = 7802

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 86
$ws.Range("L3").Value = 97
$ws.Range("L6").Value = 74
$ws.Range("L7").Value = 287

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 31
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 130
$ws.Range("L6").Value = 104
$ws.Range("L7").Value = 398

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 158
$ws.Range("K7").Value = 348

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 73
$ws.Range("L7").Value = 222

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 63
$ws.Range("L4").Value = 24
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 247

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 42
$ws.Range("L7").Value = 112

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L5").Value = 11
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 216

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 62
$ws.Range("L7").Value = 197

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 82
$ws.Range("L4").Value = 21
$ws.Range("L6").Value = 71
$ws.Range("L7").Value = 252

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 25
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 47
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 29
$ws.Range("L4").Value = 7

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 123
$ws.Range("L3").Value = 165
$ws.Range("L6").Value = 82
$ws.Range("L7").Value = 411

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 10
